# ==========================================================================
# Applies the "additional scraping" edit:
#   1. Insert a new "Player Info" sheet at the front with player bio data.
#   2. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace
#      the full scorecard URL with just the trailing MatchCode number.
#   3. On "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE treatment.
#   4. Append a new "ODI Batting Extra" sheet with extra per-match batting
#      stats (batting position, boundary counts, % of total runs, MOTM).
# ==========================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# --------------------------------------------------------------------------
# 1. "Player Info" sheet — inserted before the current first sheet.
# --------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $playerInfoHeaders[$c - 1]
}
Set-HeaderStyle $playerInfo.Range("A1:D1")

$playerInfo.Range("A2:A2").NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "3830"
$playerInfo.Cells.Item(2, 2).Value = "Mitchell Aaron Starc"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Left Arm Fast"

# --------------------------------------------------------------------------
# 2. "ODI Batting" — MATCH_CARD_LINK (col D) becomes MATCH_CODE.
# --------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastBattingRow = $batting.Cells.Item(1, 4).End(-4121).Row  # xlDown
$batting.Range("D2:D$lastBattingRow").NumberFormat = "@"
for ($r = 2; $r -le $lastBattingRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val) {
        $idx = $val.IndexOf("MatchCode=")
        if ($idx -ge 0) {
            $code = $val.Substring($idx + 10)
            $cell.Value = $code
        }
    }
}

# --------------------------------------------------------------------------
# 3. "ODI Bowling" — MATCH_CARD_LINK (col B) becomes MATCH_CODE.
# --------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastBowlingRow = $bowling.Cells.Item(1, 2).End(-4121).Row  # xlDown
$bowling.Range("B2:B$lastBowlingRow").NumberFormat = "@"
for ($r = 2; $r -le $lastBowlingRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val) {
        $idx = $val.IndexOf("MatchCode=")
        if ($idx -ge 0) {
            $code = $val.Substring($idx + 10)
            $cell.Value = $code
        }
    }
}

# --------------------------------------------------------------------------
# 4. "ODI Batting Extra" — appended after the last sheet.
# --------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $extra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
Set-HeaderStyle $extra.Range("A1:F1")

$extraLastRow = 21
$extra.Range("A2:A$extraLastRow").NumberFormat = "@"
$extra.Range("C2:C$extraLastRow").NumberFormat = "@"
$extra.Range("D2:D$extraLastRow").NumberFormat = "@"
$extra.Range("E2:E$extraLastRow").NumberFormat = "@"

# MATCH_CODE, BATTING_POSITION (numeric or $null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4423", 9,    "1", "0", "3.49%",  "NO"),
    @("4429", 9,    "1", "1", "6.46%",  "NO"),
    @("4430", 9,    "0", "0", "",       "NO"),
    @("4431", 9,    "1", "1", "3.61%",  "NO"),
    @("4435", $null, "", "", "",        "NO"),
    @("4436", $null, "", "", "",        "NO"),
    @("4483", 8,    "1", "0", "3.17%",  "YES"),
    @("4484", 8,    "2", "0", "10.16%", "NO"),
    @("4486", $null, "", "", "",        "NO"),
    @("4644", 9,    "",  "", "",        "NO"),
    @("4645", 9,    "",  "", "",        "YES"),
    @("4646", 9,    "0", "0", "1.42%",  "NO"),
    @("4647", $null, "", "", "",        "NO"),
    @("4648", 9,    "2", "1", "19.49%", "YES"),
    @("4649", 9,    "",  "", "",        "NO"),
    @("4660", $null, "", "", "",        "NO"),
    @("4663", $null, "", "", "",        "NO"),
    @("4725", 10,   "1", "0", "2.13%",  "NO"),
    @("4728", 9,    "",  "", "",        "YES"),
    @("4732", 10,   "0", "1", "3.72%",  "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    $extra.Cells.Item($r, 3).Value = $row[2]
    $extra.Cells.Item($r, 4).Value = $row[3]
    $extra.Cells.Item($r, 5).Value = $row[4]
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Keep the first sheet ("Player Info") the active/selected tab, matching the
# unchanged activeTab="0" in the workbook view.
$playerInfo.Activate() | Out-Null
$playerInfo.Range("A1").Select() | Out-Null
